# Update the "Marking"/"Total" row of the marksheet:
#  - B11 (Marking -> Right count):  3   -> 5
#  - B12 (Total   -> Right count): 72   -> 120
#  - E12 (Total   -> Correct/Total text): "72/84" -> "120/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
